# Apply the user-integration data model change:
# Add two new columns "Beobachtungen der Teilnehmenden" (R) and "Personas" (S)
# to the "German" worksheet, with data values "Teilnehmende Beobachtungen" /
# "Personas" filled in for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("German")

# New header cells
$ws.Range("R1").Value = "Beobachtungen der Teilnehmenden"
$ws.Range("S1").Value = "Personas"

# New data cells for each of the three existing data rows (2-4)
$ws.Range("R2").Value = "Teilnehmende Beobachtungen"
$ws.Range("S2").Value = "Personas"

$ws.Range("R3").Value = "Teilnehmende Beobachtungen"
$ws.Range("S3").Value = "Personas"

$ws.Range("R4").Value = "Teilnehmende Beobachtungen"
$ws.Range("S4").Value = "Personas"

# Leave the selection on the last edited cell, matching the editor's end state
$ws.Range("S2").Select()
